$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

$ws.Range("I4:I11").NumberFormat = "@"
$ws.Range("N4:N11").NumberFormat = "@"

$ws.Cells.Item(4,1).Value2 = 3
$ws.Cells.Item(4,2).Value2 = "A"
$ws.Cells.Item(4,3).Value2 = "Facility 4"
$ws.Cells.Item(4,4).Value2 = "C"
$ws.Cells.Item(4,5).Value2 = 423
$ws.Cells.Item(4,6).Value2 = 253800
$ws.Cells.Item(4,7).Value2 = "C"
$ws.Cells.Item(4,8).Value2 = 60
$ws.Cells.Item(4,9).Value2 = "4%"
$ws.Cells.Item(4,10).Value2 = 57.59999999999999
$ws.Cells.Item(4,11).Value2 = 34560
$ws.Cells.Item(4,12).Value2 = 600
$ws.Cells.Item(4,13).Value2 = 219240
$ws.Cells.Item(4,14).Value2 = "7%"
$ws.Cells.Item(4,15).Value2 = 2419.2
$ws.Cells.Item(5,1).Value2 = 4
$ws.Cells.Item(5,2).Value2 = "A"
$ws.Cells.Item(5,3).Value2 = "Facility 4"
$ws.Cells.Item(5,4).Value2 = "C"
$ws.Cells.Item(5,5).Value2 = 453
$ws.Cells.Item(5,6).Value2 = 2568510
$ws.Cells.Item(5,7).Value2 = "C"
$ws.Cells.Item(5,8).Value2 = 24
$ws.Cells.Item(5,9).Value2 = "4%"
$ws.Cells.Item(5,10).Value2 = 23.04
$ws.Cells.Item(5,11).Value2 = 130636.8
$ws.Cells.Item(5,12).Value2 = 5670
$ws.Cells.Item(5,13).Value2 = 2437873.2
$ws.Cells.Item(5,14).Value2 = "7%"
$ws.Cells.Item(5,15).Value2 = 9144.576000000001
$ws.Cells.Item(6,1).Value2 = 5
$ws.Cells.Item(6,2).Value2 = "A"
$ws.Cells.Item(6,3).Value2 = "Facility 5"
$ws.Cells.Item(6,4).Value2 = "C"
$ws.Cells.Item(6,5).Value2 = 342
$ws.Cells.Item(6,6).Value2 = 15390
$ws.Cells.Item(6,7).Value2 = "C"
$ws.Cells.Item(6,8).Value2 = 24
$ws.Cells.Item(6,9).Value2 = "4%"
$ws.Cells.Item(6,10).Value2 = 23.04
$ws.Cells.Item(6,11).Value2 = 1036.8
$ws.Cells.Item(6,12).Value2 = 45
$ws.Cells.Item(6,13).Value2 = 14353.2
$ws.Cells.Item(6,14).Value2 = "7%"
$ws.Cells.Item(6,15).Value2 = 72.57600000000001
$ws.Cells.Item(7,1).Value2 = 6
$ws.Cells.Item(7,2).Value2 = "A"
$ws.Cells.Item(7,3).Value2 = "Facility 9"
$ws.Cells.Item(7,4).Value2 = "C"
$ws.Cells.Item(7,5).Value2 = 653
$ws.Cells.Item(7,6).Value2 = 158026
$ws.Cells.Item(7,7).Value2 = "B"
$ws.Cells.Item(7,8).Value2 = 24
$ws.Cells.Item(7,9).Value2 = "0%"
$ws.Cells.Item(7,10).Value2 = 24
$ws.Cells.Item(7,11).Value2 = 5808
$ws.Cells.Item(7,12).Value2 = 242
$ws.Cells.Item(7,13).Value2 = 152218
$ws.Cells.Item(7,14).Value2 = "0%"
$ws.Cells.Item(7,15).Value2 = 0
$ws.Cells.Item(8,1).Value2 = 7
$ws.Cells.Item(8,2).Value2 = "A"
$ws.Cells.Item(8,3).Value2 = "Facility 9"
$ws.Cells.Item(8,4).Value2 = "C"
$ws.Cells.Item(8,5).Value2 = 432
$ws.Cells.Item(8,6).Value2 = 286848
$ws.Cells.Item(8,7).Value2 = "A"
$ws.Cells.Item(8,8).Value2 = 23
$ws.Cells.Item(8,9).Value2 = "1%"
$ws.Cells.Item(8,10).Value2 = 22.77
$ws.Cells.Item(8,11).Value2 = 15119.28
$ws.Cells.Item(8,12).Value2 = 664
$ws.Cells.Item(8,13).Value2 = 271728.72
$ws.Cells.Item(8,14).Value2 = "10%"
$ws.Cells.Item(8,15).Value2 = 1511.928
$ws.Cells.Item(9,1).Value2 = 8
$ws.Cells.Item(9,2).Value2 = "A"
$ws.Cells.Item(9,3).Value2 = "Facility 9"
$ws.Cells.Item(9,4).Value2 = "C"
$ws.Cells.Item(9,5).Value2 = 456
$ws.Cells.Item(9,6).Value2 = 10944
$ws.Cells.Item(9,7).Value2 = "C"
$ws.Cells.Item(9,8).Value2 = 24
$ws.Cells.Item(9,9).Value2 = "4%"
$ws.Cells.Item(9,10).Value2 = 23.04
$ws.Cells.Item(9,11).Value2 = 552.96
$ws.Cells.Item(9,12).Value2 = 24
$ws.Cells.Item(9,13).Value2 = 10391.04
$ws.Cells.Item(9,14).Value2 = "7%"
$ws.Cells.Item(9,15).Value2 = 38.70720000000001
$ws.Cells.Item(10,1).Value2 = 9
$ws.Cells.Item(10,2).Value2 = "A"
$ws.Cells.Item(10,3).Value2 = "Facility 10"
$ws.Cells.Item(10,4).Value2 = "C"
$ws.Cells.Item(10,5).Value2 = 234
$ws.Cells.Item(10,6).Value2 = 54288
$ws.Cells.Item(10,7).Value2 = "C"
$ws.Cells.Item(10,8).Value2 = 13
$ws.Cells.Item(10,9).Value2 = "4%"
$ws.Cells.Item(10,10).Value2 = 12.48
$ws.Cells.Item(10,11).Value2 = 2895.36
$ws.Cells.Item(10,12).Value2 = 232
$ws.Cells.Item(10,13).Value2 = 51392.64
$ws.Cells.Item(10,14).Value2 = "7%"
$ws.Cells.Item(10,15).Value2 = 202.6752
$ws.Cells.Item(11,1).Value2 = 10
$ws.Cells.Item(11,2).Value2 = "A"
$ws.Cells.Item(11,3).Value2 = "Facility 10"
$ws.Cells.Item(11,4).Value2 = "C"
$ws.Cells.Item(11,5).Value2 = 231
$ws.Cells.Item(11,6).Value2 = 3003
$ws.Cells.Item(11,7).Value2 = "C"
$ws.Cells.Item(11,8).Value2 = 75
$ws.Cells.Item(11,9).Value2 = "4%"
$ws.Cells.Item(11,10).Value2 = 72
$ws.Cells.Item(11,11).Value2 = 936
$ws.Cells.Item(11,12).Value2 = 13
$ws.Cells.Item(11,13).Value2 = 2067
$ws.Cells.Item(11,14).Value2 = "7%"
$ws.Cells.Item(11,15).Value2 = 65.52000000000001

$ws.Range("I4:I11").ClearFormats()
$ws.Range("N4:N11").ClearFormats()

$ws.Rows.Item(12).Delete()

$wsLP = $wb.Worksheets.Item("LP Model")
$lpText = @"
\* Sourcing_with_MultiTier_Rebates_Discounts *\
Minimize
OBJ: S_A + S_B + S_C - rebate_A - rebate_B - rebate_C
Subject To
ActiveLink_A_1: x_A_1 - 1000000000 z_A <= 0
ActiveLink_A_10: x_A_10 - 1000000000 z_A <= 0
ActiveLink_A_2: x_A_2 - 1000000000 z_A <= 0
ActiveLink_A_3: x_A_3 - 1000000000 z_A <= 0
ActiveLink_A_4: x_A_4 - 1000000000 z_A <= 0
ActiveLink_A_5: x_A_5 - 1000000000 z_A <= 0
ActiveLink_A_6: x_A_6 - 1000000000 z_A <= 0
ActiveLink_A_7: x_A_7 - 1000000000 z_A <= 0
ActiveLink_A_8: x_A_8 - 1000000000 z_A <= 0
ActiveLink_A_9: x_A_9 - 1000000000 z_A <= 0
ActiveLink_B_1: x_B_1 - 1000000000 z_B <= 0
ActiveLink_B_10: x_B_10 - 1000000000 z_B <= 0
ActiveLink_B_2: x_B_2 - 1000000000 z_B <= 0
ActiveLink_B_3: x_B_3 - 1000000000 z_B <= 0
ActiveLink_B_4: x_B_4 - 1000000000 z_B <= 0
ActiveLink_B_5: x_B_5 - 1000000000 z_B <= 0
ActiveLink_B_6: x_B_6 - 1000000000 z_B <= 0
ActiveLink_B_7: x_B_7 - 1000000000 z_B <= 0
ActiveLink_B_8: x_B_8 - 1000000000 z_B <= 0
ActiveLink_B_9: x_B_9 - 1000000000 z_B <= 0
ActiveLink_C_1: x_C_1 - 1000000000 z_C <= 0
ActiveLink_C_10: x_C_10 - 1000000000 z_C <= 0
ActiveLink_C_2: x_C_2 - 1000000000 z_C <= 0
ActiveLink_C_3: x_C_3 - 1000000000 z_C <= 0
ActiveLink_C_4: x_C_4 - 1000000000 z_C <= 0
ActiveLink_C_5: x_C_5 - 1000000000 z_C <= 0
ActiveLink_C_6: x_C_6 - 1000000000 z_C <= 0
ActiveLink_C_7: x_C_7 - 1000000000 z_C <= 0
ActiveLink_C_8: x_C_8 - 1000000000 z_C <= 0
ActiveLink_C_9: x_C_9 - 1000000000 z_C <= 0
BaseSpend_A: S0_A - 50 x_A_1 - 64 x_A_10 - 70 x_A_2 - 55 x_A_3 - 23 x_A_4
 - 54 x_A_5 - 42 x_A_6 - 23 x_A_7 - 75 x_A_8 - 97 x_A_9 = 0
BaseSpend_B: S0_B - 60 x_B_1 - 13 x_B_10 - 80 x_B_2 - 65 x_B_3 - 75 x_B_4
 - 34 x_B_5 - 24 x_B_6 - 53 x_B_7 - 13 x_B_8 - 56 x_B_9 = 0
BaseSpend_C: S0_C - 55 x_C_1 - 75 x_C_10 - 75 x_C_2 - 60 x_C_3 - 24 x_C_4
 - 24 x_C_5 - 64 x_C_6 - 86 x_C_7 - 24 x_C_8 - 13 x_C_9 = 0
Capacity_B_Bid_ID_1: x_B_1 <= 100000000
Capacity_B_Bid_ID_10: x_B_10 <= 100000000
Capacity_B_Bid_ID_2: x_B_2 <= 100000000
Capacity_B_Bid_ID_3: x_B_3 <= 100000000
Capacity_B_Bid_ID_4: x_B_4 <= 100000000
Capacity_B_Bid_ID_5: x_B_5 <= 100000000
Capacity_B_Bid_ID_6: x_B_6 <= 100000000
Capacity_B_Bid_ID_7: x_B_7 <= 100000000
Capacity_B_Bid_ID_8: x_B_8 <= 100000000
Capacity_B_Bid_ID_9: x_B_9 <= 100000000
Capacity_C_Bid_ID_1: x_C_1 <= 100000000
Capacity_C_Bid_ID_10: x_C_10 <= 100000000
Capacity_C_Bid_ID_2: x_C_2 <= 100000000
Capacity_C_Bid_ID_3: x_C_3 <= 100000000
Capacity_C_Bid_ID_4: x_C_4 <= 100000000
Capacity_C_Bid_ID_5: x_C_5 <= 100000000
Capacity_C_Bid_ID_6: x_C_6 <= 100000000
Capacity_C_Bid_ID_7: x_C_7 <= 100000000
Capacity_C_Bid_ID_8: x_C_8 <= 100000000
Capacity_C_Bid_ID_9: x_C_9 <= 100000000
Demand_1: x_A_1 + x_B_1 + x_C_1 = 700
Demand_10: x_A_10 + x_B_10 + x_C_10 = 13
Demand_2: x_A_2 + x_B_2 + x_C_2 = 9000
Demand_3: x_A_3 + x_B_3 + x_C_3 = 600
Demand_4: x_A_4 + x_B_4 + x_C_4 = 5670
Demand_5: x_A_5 + x_B_5 + x_C_5 = 45
Demand_6: x_A_6 + x_B_6 + x_C_6 = 242
Demand_7: x_A_7 + x_B_7 + x_C_7 = 664
Demand_8: x_A_8 + x_B_8 + x_C_8 = 24
Demand_9: x_A_9 + x_B_9 + x_C_9 = 232
DiscountTierLower_A_0: d_A - 19400000000 z_discount_A_0 >= -19400000000
DiscountTierLower_A_1: - 0.01 S0_A + d_A - 19400000000 z_discount_A_1
 >= -19400000000
DiscountTierLower_B_0: d_B - 97000000000 z_discount_B_0 >= -97000000000
DiscountTierLower_B_1: - 0.03 S0_B + d_B - 97000000000 z_discount_B_1
 >= -97000000000
DiscountTierLower_C_0: d_C - 97000000000 z_discount_C_0 >= -97000000000
DiscountTierLower_C_1: - 0.04 S0_C + d_C - 97000000000 z_discount_C_1
 >= -97000000000
DiscountTierMax_A_0: x_A_1 + x_A_10 + x_A_2 + x_A_3 + x_A_4 + x_A_5 + x_A_6
 + x_A_7 + x_A_8 + x_A_9 + 19400000000 z_discount_A_0 <= 19400001000
DiscountTierMax_B_0: x_B_1 + x_B_10 + x_B_2 + x_B_3 + x_B_4 + x_B_5 + x_B_6
 + x_B_7 + x_B_8 + x_B_9 + 97000000000 z_discount_B_0 <= 97000000500
DiscountTierMax_C_0: x_C_1 + x_C_10 + x_C_2 + x_C_3 + x_C_4 + x_C_5 + x_C_6
 + x_C_7 + x_C_8 + x_C_9 + 97000000000 z_discount_C_0 <= 97000000500
DiscountTierMin_A_0: x_A_1 + x_A_10 + x_A_2 + x_A_3 + x_A_4 + x_A_5 + x_A_6
 + x_A_7 + x_A_8 + x_A_9 >= 0
DiscountTierMin_A_1: x_A_1 + x_A_10 + x_A_2 + x_A_3 + x_A_4 + x_A_5 + x_A_6
 + x_A_7 + x_A_8 + x_A_9 - 1000 z_discount_A_1 >= 0
DiscountTierMin_B_0: x_B_1 + x_B_10 + x_B_2 + x_B_3 + x_B_4 + x_B_5 + x_B_6
 + x_B_7 + x_B_8 + x_B_9 >= 0
DiscountTierMin_B_1: x_B_1 + x_B_10 + x_B_2 + x_B_3 + x_B_4 + x_B_5 + x_B_6
 + x_B_7 + x_B_8 + x_B_9 - 500 z_discount_B_1 >= 0
DiscountTierMin_C_0: x_C_1 + x_C_10 + x_C_2 + x_C_3 + x_C_4 + x_C_5 + x_C_6
 + x_C_7 + x_C_8 + x_C_9 >= 0
DiscountTierMin_C_1: x_C_1 + x_C_10 + x_C_2 + x_C_3 + x_C_4 + x_C_5 + x_C_6
 + x_C_7 + x_C_8 + x_C_9 - 500 z_discount_C_1 >= 0
DiscountTierSelect_A: z_discount_A_0 + z_discount_A_1 = 1
DiscountTierSelect_B: z_discount_B_0 + z_discount_B_1 = 1
DiscountTierSelect_C: z_discount_C_0 + z_discount_C_1 = 1
DiscountTierUpper_A_0: d_A + 19400000000 z_discount_A_0 <= 19400000000
DiscountTierUpper_A_1: - 0.01 S0_A + d_A + 19400000000 z_discount_A_1
 <= 19400000000
DiscountTierUpper_B_0: d_B + 97000000000 z_discount_B_0 <= 97000000000
DiscountTierUpper_B_1: - 0.03 S0_B + d_B + 97000000000 z_discount_B_1
 <= 97000000000
DiscountTierUpper_C_0: d_C + 97000000000 z_discount_C_0 <= 97000000000
DiscountTierUpper_C_1: - 0.04 S0_C + d_C + 97000000000 z_discount_C_1
 <= 97000000000
EffectiveSpend_A: - S0_A + S_A + d_A = 0
EffectiveSpend_B: - S0_B + S_B + d_B = 0
EffectiveSpend_C: - S0_C + S_C + d_C = 0
MinAward_A: x_A_1 + x_A_10 + x_A_2 + x_A_3 + x_A_4 + x_A_5 + x_A_6 + x_A_7
 + x_A_8 + x_A_9 - z_A >= 0
MinAward_B: x_B_1 + x_B_10 + x_B_2 + x_B_3 + x_B_4 + x_B_5 + x_B_6 + x_B_7
 + x_B_8 + x_B_9 - z_B >= 0
MinAward_C: x_C_1 + x_C_10 + x_C_2 + x_C_3 + x_C_4 + x_C_5 + x_C_6 + x_C_7
 + x_C_8 + x_C_9 - z_C >= 0
RebateTierLower_A_0: rebate_A - 19400000000 y_rebate_A_0 >= -19400000000
RebateTierLower_A_1: - 0.1 S_A + rebate_A - 19400000000 y_rebate_A_1
 >= -19400000000
RebateTierLower_B_0: rebate_B - 97000000000 y_rebate_B_0 >= -97000000000
RebateTierLower_B_1: - 0.05 S_B + rebate_B - 97000000000 y_rebate_B_1
 >= -97000000000
RebateTierLower_C_0: rebate_C - 97000000000 y_rebate_C_0 >= -97000000000
RebateTierLower_C_1: - 0.07 S_C + rebate_C - 97000000000 y_rebate_C_1
 >= -97000000000
RebateTierMax_A_0: x_A_1 + x_A_10 + x_A_2 + x_A_3 + x_A_4 + x_A_5 + x_A_6
 + x_A_7 + x_A_8 + x_A_9 + 19400000000 y_rebate_A_0 <= 19400000500
RebateTierMax_B_0: x_B_1 + x_B_10 + x_B_2 + x_B_3 + x_B_4 + x_B_5 + x_B_6
 + x_B_7 + x_B_8 + x_B_9 + 97000000000 y_rebate_B_0 <= 97000000500
RebateTierMax_C_0: x_C_1 + x_C_10 + x_C_2 + x_C_3 + x_C_4 + x_C_5 + x_C_6
 + x_C_7 + x_C_8 + x_C_9 + 97000000000 y_rebate_C_0 <= 97000000700
RebateTierMin_A_0: x_A_1 + x_A_10 + x_A_2 + x_A_3 + x_A_4 + x_A_5 + x_A_6
 + x_A_7 + x_A_8 + x_A_9 >= 0
RebateTierMin_A_1: x_A_1 + x_A_10 + x_A_2 + x_A_3 + x_A_4 + x_A_5 + x_A_6
 + x_A_7 + x_A_8 + x_A_9 - 500 y_rebate_A_1 >= 0
RebateTierMin_B_0: x_B_1 + x_B_10 + x_B_2 + x_B_3 + x_B_4 + x_B_5 + x_B_6
 + x_B_7 + x_B_8 + x_B_9 >= 0
RebateTierMin_B_1: x_B_1 + x_B_10 + x_B_2 + x_B_3 + x_B_4 + x_B_5 + x_B_6
 + x_B_7 + x_B_8 + x_B_9 - 500 y_rebate_B_1 >= 0
RebateTierMin_C_0: x_C_1 + x_C_10 + x_C_2 + x_C_3 + x_C_4 + x_C_5 + x_C_6
 + x_C_7 + x_C_8 + x_C_9 >= 0
RebateTierMin_C_1: x_C_1 + x_C_10 + x_C_2 + x_C_3 + x_C_4 + x_C_5 + x_C_6
 + x_C_7 + x_C_8 + x_C_9 - 700 y_rebate_C_1 >= 0
RebateTierSelect_A: y_rebate_A_0 + y_rebate_A_1 = 1
RebateTierSelect_B: y_rebate_B_0 + y_rebate_B_1 = 1
RebateTierSelect_C: y_rebate_C_0 + y_rebate_C_1 = 1
RebateTierUpper_A_0: rebate_A + 19400000000 y_rebate_A_0 <= 19400000000
RebateTierUpper_A_1: - 0.1 S_A + rebate_A + 19400000000 y_rebate_A_1
 <= 19400000000
RebateTierUpper_B_0: rebate_B + 97000000000 y_rebate_B_0 <= 97000000000
RebateTierUpper_B_1: - 0.05 S_B + rebate_B + 97000000000 y_rebate_B_1
 <= 97000000000
RebateTierUpper_C_0: rebate_C + 97000000000 y_rebate_C_0 <= 97000000000
RebateTierUpper_C_1: - 0.07 S_C + rebate_C + 97000000000 y_rebate_C_1
 <= 97000000000
TransitionLower_10_A: x_A_10 >= 0
TransitionLower_10_B: x_B_10 >= 0
TransitionLower_1_B: x_B_1 >= 0
TransitionLower_1_C: x_C_1 >= 0
TransitionLower_2_A: x_A_2 >= 0
TransitionLower_2_C: x_C_2 >= 0
TransitionLower_3_A: x_A_3 >= 0
TransitionLower_3_B: x_B_3 >= 0
TransitionLower_4_A: x_A_4 >= 0
TransitionLower_4_B: x_B_4 >= 0
TransitionLower_5_A: x_A_5 >= 0
TransitionLower_5_B: x_B_5 >= 0
TransitionLower_6_A: x_A_6 >= 0
TransitionLower_6_B: x_B_6 >= 0
TransitionLower_7_A: x_A_7 >= 0
TransitionLower_7_B: x_B_7 >= 0
TransitionLower_8_A: x_A_8 >= 0
TransitionLower_8_B: x_B_8 >= 0
TransitionLower_9_A: x_A_9 >= 0
TransitionLower_9_B: x_B_9 >= 0
Transition_10_A: - 13 T_10_A + x_A_10 <= 0
Transition_10_B: - 13 T_10_B + x_B_10 <= 0
Transition_1_B: - 700 T_1_B + x_B_1 <= 0
Transition_1_C: - 700 T_1_C + x_C_1 <= 0
Transition_2_A: - 9000 T_2_A + x_A_2 <= 0
Transition_2_C: - 9000 T_2_C + x_C_2 <= 0
Transition_3_A: - 600 T_3_A + x_A_3 <= 0
Transition_3_B: - 600 T_3_B + x_B_3 <= 0
Transition_4_A: - 5670 T_4_A + x_A_4 <= 0
Transition_4_B: - 5670 T_4_B + x_B_4 <= 0
Transition_5_A: - 45 T_5_A + x_A_5 <= 0
Transition_5_B: - 45 T_5_B + x_B_5 <= 0
Transition_6_A: - 242 T_6_A + x_A_6 <= 0
Transition_6_B: - 242 T_6_B + x_B_6 <= 0
Transition_7_A: - 664 T_7_A + x_A_7 <= 0
Transition_7_B: - 664 T_7_B + x_B_7 <= 0
Transition_8_A: - 24 T_8_A + x_A_8 <= 0
Transition_8_B: - 24 T_8_B + x_B_8 <= 0
Transition_9_A: - 232 T_9_A + x_A_9 <= 0
Transition_9_B: - 232 T_9_B + x_B_9 <= 0
Transitions_0_UB: T_10_A + T_10_B + T_3_A + T_3_B + T_4_A + T_4_B + T_5_A
 + T_5_B + T_6_A + T_6_B + T_7_A + T_7_B + T_8_A + T_8_B + T_9_A + T_9_B <= 2
Volume_A: V_A - x_A_1 - x_A_10 - x_A_2 - x_A_3 - x_A_4 - x_A_5 - x_A_6 - x_A_7
 - x_A_8 - x_A_9 = 0
Volume_B: V_B - x_B_1 - x_B_10 - x_B_2 - x_B_3 - x_B_4 - x_B_5 - x_B_6 - x_B_7
 - x_B_8 - x_B_9 = 0
Volume_C: V_C - x_C_1 - x_C_10 - x_C_2 - x_C_3 - x_C_4 - x_C_5 - x_C_6 - x_C_7
 - x_C_8 - x_C_9 = 0
Binaries
T_10_A
T_10_B
T_1_B
T_1_C
T_2_A
T_2_C
T_3_A
T_3_B
T_4_A
T_4_B
T_5_A
T_5_B
T_6_A
T_6_B
T_7_A
T_7_B
T_8_A
T_8_B
T_9_A
T_9_B
y_rebate_A_0
y_rebate_A_1
y_rebate_B_0
y_rebate_B_1
y_rebate_C_0
y_rebate_C_1
z_A
z_B
z_C
z_discount_A_0
z_discount_A_1
z_discount_B_0
z_discount_B_1
z_discount_C_0
z_discount_C_1
End

"@
$wsLP.Range("A2").Value2 = $lpText
